$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for new columns I (I0) and J (IF), styled like existing header row (H1 uses style index 1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Fill data rows 2-37: I = 1 (constant), J = same value as column H
for ($r = 2; $r -le 37; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
